# Add new columns I (I0) and J (IF) to the worksheet, matching the
# header style used by the existing H1 ("IP") header cell, and fill in
# the per-row values for rows 2-82.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Match the formatting of the existing header cells (bold, centered,
# top-aligned, thin border) used by B1:H1.
$ws.Range("I1:J1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("I1:J1").VerticalAlignment = -4160    # xlTop
$ws.Range("I1:J1").Borders.LineStyle = 1

# --- Data rows (rows 2-82) ---
$colI = @(7,8,7,8,8,8,8,8,8,7,9,8,8,8,8,8,8,9,8,8,7,8,9,8,11,7,8,8,7,8,7,8,8,10,8,9,9,9,9,5,7,2,9,9,8,7,9,9,9,9,7,6,8,5,6,7,7,7,9,9,9,6,6,7,5,6,7,7,5,8,7,7,9,5,6,9,4,9,9,5,6)
$colJ = @(8,8,8,8,8,8,8,8,8,8,9,8,8,8,8,8,8,9,8,8,8,8,9,8,11,8,8,8,7,8,8,8,8,10,9,9,9,9,9,6,7,3,9,9,8,8,9,9,9,9,8,7,9,7,8,7,8,7,9,9,9,7,6,8,5,6,8,7,5,8,7,7,9,6,6,9,4,9,9,5,6)

for ($i = 0; $i -lt $colI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $colI[$i]
    $ws.Cells.Item($row, 10).Value = $colJ[$i]
}
